$d = $word.ActiveDocument

# Update the date line
$d.Content.Find.Execute("2024-03-01 Friday", $true, $false, $false, $false, $false, $true, 1, $false, "2024-03-02 Saturday", 2) | Out-Null

# Update the table of math problems by cell reference to avoid any ambiguous text matches
$tbl = $d.Tables.Item(1)

$tbl.Cell(1, 1).Range.Text = "5+78="
$tbl.Cell(1, 2).Range.Text = "44+27="
$tbl.Cell(1, 3).Range.Text = "81-48="
$tbl.Cell(1, 4).Range.Text = "39-39="
$tbl.Cell(1, 5).Range.Text = "22-6="
$tbl.Cell(2, 1).Range.Text = "25+50="
$tbl.Cell(2, 2).Range.Text = "23+75="
$tbl.Cell(2, 3).Range.Text = "38+14="
$tbl.Cell(2, 4).Range.Text = "53-34="
$tbl.Cell(2, 5).Range.Text = "31+31="
$tbl.Cell(3, 1).Range.Text = "83-43="
$tbl.Cell(3, 2).Range.Text = "98-51="
$tbl.Cell(3, 3).Range.Text = "55-23="
$tbl.Cell(3, 4).Range.Text = "5+78="
$tbl.Cell(3, 5).Range.Text = "10+28="
$tbl.Cell(4, 1).Range.Text = "44+46="
$tbl.Cell(4, 2).Range.Text = "35-3="
$tbl.Cell(4, 3).Range.Text = "99-42="
$tbl.Cell(4, 4).Range.Text = "60-6="
$tbl.Cell(4, 5).Range.Text = "39+20="
$tbl.Cell(5, 1).Range.Text = "17+38="
$tbl.Cell(5, 2).Range.Text = "18+69="
$tbl.Cell(5, 3).Range.Text = "11+42="
$tbl.Cell(5, 4).Range.Text = "26-21="
$tbl.Cell(5, 5).Range.Text = "67-44="
$tbl.Cell(6, 1).Range.Text = "25+0="
$tbl.Cell(6, 2).Range.Text = "71+4="
$tbl.Cell(6, 3).Range.Text = "64-53="
$tbl.Cell(6, 4).Range.Text = "55-44="
$tbl.Cell(6, 5).Range.Text = "19-7="
$tbl.Cell(7, 1).Range.Text = "66-30="
$tbl.Cell(7, 2).Range.Text = "8+77="
$tbl.Cell(7, 3).Range.Text = "19+47="
$tbl.Cell(7, 4).Range.Text = "15+12="
$tbl.Cell(7, 5).Range.Text = "69+27="
$tbl.Cell(8, 1).Range.Text = "6-2="
$tbl.Cell(8, 2).Range.Text = "33-11="
$tbl.Cell(8, 3).Range.Text = "20+70="
$tbl.Cell(8, 4).Range.Text = "6+82="
$tbl.Cell(8, 5).Range.Text = "60-42="
$tbl.Cell(9, 1).Range.Text = "38+58="
$tbl.Cell(9, 2).Range.Text = "40+9="
$tbl.Cell(9, 3).Range.Text = "36+1="
$tbl.Cell(9, 4).Range.Text = "10+82="
$tbl.Cell(9, 5).Range.Text = "32+59="
$tbl.Cell(10, 1).Range.Text = "73-17="
$tbl.Cell(10, 2).Range.Text = "48+15="
$tbl.Cell(10, 3).Range.Text = "15+69="
$tbl.Cell(10, 4).Range.Text = "16+1="
$tbl.Cell(10, 5).Range.Text = "36+11="
$tbl.Cell(11, 1).Range.Text = "14+23="
$tbl.Cell(11, 2).Range.Text = "0+13="
$tbl.Cell(11, 3).Range.Text = "23+52="
$tbl.Cell(11, 4).Range.Text = "11+11="
$tbl.Cell(11, 5).Range.Text = "0+27="
$tbl.Cell(12, 1).Range.Text = "38+7="
$tbl.Cell(12, 2).Range.Text = "59+12="
$tbl.Cell(12, 3).Range.Text = "75-48="
$tbl.Cell(12, 4).Range.Text = "9+70="
$tbl.Cell(12, 5).Range.Text = "68-24="
$tbl.Cell(13, 1).Range.Text = "62+6="
$tbl.Cell(13, 2).Range.Text = "87+0="
$tbl.Cell(13, 3).Range.Text = "61-20="
$tbl.Cell(13, 4).Range.Text = "25-18="
$tbl.Cell(13, 5).Range.Text = "60-0="
$tbl.Cell(14, 1).Range.Text = "28-20="
$tbl.Cell(14, 2).Range.Text = "32+38="
$tbl.Cell(14, 3).Range.Text = "94-57="
$tbl.Cell(14, 4).Range.Text = "65+11="
$tbl.Cell(14, 5).Range.Text = "29+14="
$tbl.Cell(15, 1).Range.Text = "49-39="
$tbl.Cell(15, 2).Range.Text = "51+14="
$tbl.Cell(15, 3).Range.Text = "12+35="
$tbl.Cell(15, 4).Range.Text = "3+67="
$tbl.Cell(15, 5).Range.Text = "78-22="
$tbl.Cell(16, 1).Range.Text = "42+22="
$tbl.Cell(16, 2).Range.Text = "20+73="
$tbl.Cell(16, 3).Range.Text = "81-75="
$tbl.Cell(16, 4).Range.Text = "15+13="
$tbl.Cell(16, 5).Range.Text = "61+8="
$tbl.Cell(17, 1).Range.Text = "80-49="
$tbl.Cell(17, 2).Range.Text = "2+53="
$tbl.Cell(17, 3).Range.Text = "10+79="
$tbl.Cell(17, 4).Range.Text = "94-29="
$tbl.Cell(17, 5).Range.Text = "49+22="
$tbl.Cell(18, 1).Range.Text = "69+5="
$tbl.Cell(18, 2).Range.Text = "51+31="
$tbl.Cell(18, 3).Range.Text = "36+17="
$tbl.Cell(18, 4).Range.Text = "93-15="
$tbl.Cell(18, 5).Range.Text = "50-42="
$tbl.Cell(19, 1).Range.Text = "43-15="
$tbl.Cell(19, 2).Range.Text = "38+13="
$tbl.Cell(19, 3).Range.Text = "60-35="
$tbl.Cell(19, 4).Range.Text = "18+54="
$tbl.Cell(19, 5).Range.Text = "35-16="
$tbl.Cell(20, 1).Range.Text = "84-58="
$tbl.Cell(20, 2).Range.Text = "64-47="
$tbl.Cell(20, 3).Range.Text = "46+0="
$tbl.Cell(20, 4).Range.Text = "21+41="
$tbl.Cell(20, 5).Range.Text = "71-44="
